$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 410 (pushes existing rows 410:509 down to 411:510)
$ws.Rows(410).Insert()

# Populate the newly inserted row 410 with the new data record
$ws.Range("A410").Value = 10
$ws.Range("B410").Value = "Vega Modelo de Temuco"
$ws.Range("C410").Value = "La Araucanía"
$ws.Range("D410").Value = 45135
$ws.Range("E410").Value = 9
$ws.Range("F410").Value = 100114013
$ws.Range("G410").Value = "Zanahoria"
$ws.Range("H410").Value = "Sin especificar"
$ws.Range("I410").Value = "Primera"
$ws.Range("J410").Value = 250
$ws.Range("K410").Value = 5000
$ws.Range("L410").Value = 5000
$ws.Range("M410").Value = 5000
$ws.Range("N410").Value = "`$/saco 25 kilos"
$ws.Range("O410").Value = "Región de La Araucanía"
$ws.Range("P410").Value = 200
$ws.Range("Q410").Value = 25
$ws.Range("R410").Value = "Hortaliza"
